$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '21.771.85'
$ws.Range('D2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.574.22'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +6.98%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9867'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '286.23'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('E6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3692'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.26%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3275'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +7.06%  '
$ws.Range('E8').ClearFormats()

$ws.Range('B9').Value = 'Polygon'

$ws.Range('C9').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.142'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.97%  '
$ws.Range('E9').ClearFormats()

$ws.Range('B10').Value = 'OKB'

$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.72'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07044'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +6.27%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9984'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.06'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +10.92%  '
$ws.Range('E13').ClearFormats()

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.73%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.516'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.42%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001067'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.77%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9861'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.569.55'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.58%  '
$ws.Range('E18').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.07%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '74.70'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.09%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.20'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +11.60%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.876'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +7.77%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.62'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.93%  '
$ws.Range('E23').ClearFormats()

$ws.Range('B24').Value = 'WrappedBTC'

$ws.Range('C24').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '21.784.08'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.80%  '
$ws.Range('E24').ClearFormats()

$ws.Range('B25').Value = 'Toncoin'

$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.355'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.91%  '
$ws.Range('E25').ClearFormats()

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.383'
$ws.Range('D26').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +12.25%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.93'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.84%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.43'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +7.05%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.742.65'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.93%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.11'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.57%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.105'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.76%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9139'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +12.65%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.450'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +9.66%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08204'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.601'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.02%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.139'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +9.28%  '
$ws.Range('E36').ClearFormats()

$ws.Range('B37').Value = 'Aptos'

$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.61'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.34%  '
$ws.Range('E37').ClearFormats()

$ws.Range('B38').Value = 'Hedera'

$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06074'
$ws.Range('D38').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.31%  '
$ws.Range('E38').ClearFormats()

$ws.Range('B39').Value = 'TrustWalletToken'

$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.235'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.90%  '
$ws.Range('E39').ClearFormats()

$ws.Range('B40').Value = 'FraxShare'

$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.297'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +9.32%  '
$ws.Range('E40').ClearFormats()

$ws.Range('B41').Value = 'VeChain'

$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02172'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.69%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2001'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.54%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9860'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.84%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5761'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +9.20%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.89'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.82%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.628'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5638'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +8.91%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.59'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.76%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.913'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +7.22%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06744'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.39%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.86'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +7.31%  '
$ws.Range('E51').ClearFormats()
